# drawingML/VML export: fix position of shape in case rotation is 180 degrees.
# This introduces a new elbow-connector shape (_x0000_s1035, rotation:180)
# together with its shapetype definition, placed right after the "Hardware"
# rectangle (_x0000_s1033) and before the pre-existing _x0000_t34 shapetype
# definition inside the drawing's v:group.

$d = $word.ActiveDocument

# Work on the whole document as a single Range so we can read/replace its
# full OOXML (flat-OPC package) representation, which is the only reliable
# way to touch the VML (w:pict) content that isn't exposed through the
# normal Shapes/Range text object model.
$full = $d.Range(0, $d.Content.End)
$xml = $full.WordOpenXML

# New shapetype + shape markup taken from the target diff.
$newShapeType = '<v:shapetype id="_x0000_t34" coordsize="21600,21600" o:spt="34" o:oned="t" adj="10800" path="m,l@0,0@0,21600,21600,21600e" filled="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="val #0"/></v:formulas><v:path arrowok="t" fillok="f" o:connecttype="none"/><v:handles><v:h position="#0,center"/></v:handles><o:lock v:ext="edit" shapetype="t"/></v:shapetype>'
$newShape = '<v:shape id="_x0000_s1035" type="#_x0000_t34" style="position:absolute;left:2956;top:291;width:1;height:495;rotation:180" o:connectortype="elbow" adj="-7776000,-486628,77954400"><v:stroke startarrow="block" endarrow="block"/></v:shape>'

$insertion = $newShapeType + $newShape

# Insert right after the rect (_x0000_s1033) closes and before the existing
# <v:shapetype id="_x0000_t34" ...> element. Using a regex keeps this
# resilient to the exact whitespace/indentation used between the tags.
$pattern = '(</v:rect>\s*)(<v:shapetype id="_x0000_t34")'
$newXml = $xml -replace $pattern, ('${1}' + $insertion + '${2}')

if ($newXml -eq $xml) {
    Write-Host "ERROR: pattern not found, document left unchanged"
} else {
    $full.InsertXML($newXml)
    Write-Host "OK: inserted new v:shapetype/v:shape (_x0000_s1035) before existing shapetype."
}
